$d = $word.ActiveDocument

# --- Step 1: remove the stray _GoBack bookmark near "Step 5: Your Brand Statement" ---
$step5Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="75A34747" w14:textId="4C768994" w:rsidR="00253B63" w:rsidRPr="00253B63" w:rsidRDefault="00253B63" w:rsidP="00253B63"><w:pPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00253B63"><w:rPr><w:rFonts w:ascii="Roboto" w:eastAsia="Times New Roman" w:hAnsi="Roboto" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Step 5: Your Brand Statement</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$step5Range = $d.Content
$found = $step5Range.Find.Execute("Step 5: Your Brand Statement", $false)
if ($found) {
    $step5Para = $step5Range.Paragraphs.Item(1).Range
    $step5Para.Collapse(1)
    $step5Para.InsertXML($step5Xml)
}

# --- Step 2: rewrite the brand-statement paragraph with new wording + relocated _GoBack bookmark ---
$paraXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="4DDB6C76" w14:textId="3EEE32B6" w:rsidR="00E25FB4" w:rsidRDefault="00E25FB4" w:rsidP="00E25FB4"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t>Energetic professional marketer with a</w:t></w:r><w:r w:rsidR="00C96D9A"><w:t xml:space="preserve"> passion for learning from and </w:t></w:r><w:r w:rsidR="00BC25DD"><w:t>working</w:t></w:r><w:r w:rsidR="00C96D9A"><w:t xml:space="preserve"> with </w:t></w:r><w:r w:rsidR="00F17ED9"><w:t>a wide variety of people</w:t></w:r><w:r w:rsidR="00C96D9A"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00C866DA"><w:t xml:space="preserve">has over a decade of </w:t></w:r><w:r w:rsidR="00A17619"><w:t xml:space="preserve">multi-industry </w:t></w:r><w:r w:rsidR="00C866DA"><w:t xml:space="preserve">experience in </w:t></w:r><w:r w:rsidR="00E853D7"><w:t xml:space="preserve">marketing, administrative support and service </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>seeks to aide small</w:t></w:r><w:r w:rsidR="00674B75"><w:t xml:space="preserve"> to medium sized</w:t></w:r><w:r><w:t xml:space="preserve"> businesses </w:t></w:r><w:r><w:t xml:space="preserve">utilizing data engineering to improve and enhance strategic </w:t></w:r><w:r w:rsidR="00A17619"><w:t xml:space="preserve">marketing </w:t></w:r><w:r><w:t>plan</w:t></w:r><w:r w:rsidR="00960272"><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00960272"><w:t xml:space="preserve">to </w:t></w:r><w:r w:rsidR="00A17619"><w:t>create awareness, drive purchases, and generate loyalty among consumers/customers.</w:t></w:r><w:r w:rsidR="0003317B"><w:t xml:space="preserve">  </w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange = $d.Content
$found2 = $targetRange.Find.Execute("Energetic professional marketer with a passion for learning", $false)
if ($found2) {
    $targetPara = $targetRange.Paragraphs.Item(1).Range
    $targetPara.Collapse(1)
    $targetPara.InsertXML($paraXml)
}

Write-Output "done"
